# Add three new log rows (31-33) to the Project Log sheet: "Labyrinth",
# "Torches" and "Added RAIN AI" entries, matching the author's original
# typing order so the shared-string table is rebuilt in the same sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33 first: "Added RAIN AI" / AI patrol note ---
$ws.Range("A33").Value = "Added RAIN AI"
$ws.Range("B33").Value = 43012
$ws.Range("E33").Value = "AI patrol route added, currently will clip through wall"

# --- Row 31 next: "Labyrinth" fix note ---
$ws.Range("A31").Value = "Labyrinth"
$ws.Range("B31").Value = 42798
$ws.Range("C31").Value = 42951
$ws.Range("E31").Value = "Fixed and added more things to labyrinth, made it look better"

# --- Row 32 last: "Torches" note ---
$ws.Range("B32").Value = 42798
$ws.Range("C32").Value = 42951
$ws.Range("E32").Value = "Added torches"
$ws.Range("A32").Value = "Torches"

# Give the new date cells the same formatting (short date, style index 1)
# already used elsewhere in the sheet, by copying an existing date cell's
# format instead of assigning a brand-new number format.
$ws.Range("B25:C25").Copy()
$ws.Range("B31:C32").PasteSpecial(-4122)
$ws.Range("B25").Copy()
$ws.Range("B33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the selection where the author's last entry landed.
$ws.Range("D36").Select() | Out-Null
